$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "NO."
$ws.Range("B1").Value = "questions"
$ws.Range("C1").Value = "choice1"
$ws.Range("D1").Value = "choice2"
$ws.Range("E1").Value = "choice3"
$ws.Range("F1").Value = "choice4"
$ws.Range("G1").Value = "answer"
$ws.Range("H1").Value = "image"

# Header row uses the plain Arial style (style index 0), so make sure
# A1:C1 (which previously used the Lohit Devanagari style) switch back.
$ws.Range("A1:H1").Font.Name = "Arial"

# --- Row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "خوب هستید؟"
$ws.Range("C2").Value = "اره"
$ws.Range("D2").Value = "نه"
$ws.Range("E2").Value = "yes"
$ws.Range("F2").Value = "no"
$ws.Range("G2").Value = "yes"
$ws.Range("H2").Value = "none"

$ws.Range("A2").Font.Name = "Arial"
$ws.Range("B2:D2").Font.Name = "Lohit Devanagari"
$ws.Range("E2:H2").Font.Name = "Arial"

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "چی کارا میکنی؟"
$ws.Range("C3").Value = "بیکار"
$ws.Range("D3").Value = "علاف"
$ws.Range("E3").Value = "درس"
$ws.Range("F3").Value = "زبان"
$ws.Range("G3").Value = "درس"
$ws.Range("H3").Value = "none"

$ws.Range("A3").Font.Name = "Arial"
$ws.Range("B3:G3").Font.Name = "Lohit Devanagari"
$ws.Range("H3").Font.Name = "Arial"

# --- Selection ---
$ws.Range("H2").Select()
